# Recompute column H ("客単価" per visit) on the ABC分析_客構成 sheet so that
# it reflects spend per visit (B / (C * E)) instead of spend per unique
# customer (B / C). Column H currently holds static numeric values (no
# formulas), so we recompute each affected cell directly from columns
# B (total sales), C (customer count) and E (visit count) already on the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABC分析_客構成")

for ($r = 2; $r -le 42; $r++) {
    $cCell = $ws.Cells.Item($r, 3)   # column C
    $cVal = $cCell.Value()

    # Skip rows where C is 0 (division by zero => "inf" text, left untouched)
    if ($cVal -eq 0) {
        continue
    }

    $hCell = $ws.Cells.Item($r, 8)   # column H
    # Only rewrite cells that currently hold a numeric value
    if ($hCell.Value() -eq $null) {
        continue
    }

    $bVal = $ws.Cells.Item($r, 2).Value()  # column B
    $eVal = $ws.Cells.Item($r, 5).Value()  # column E

    $hCell.Value = $bVal / ($cVal * $eVal)
}
